# ---------------------------------------------------------------------------
# EmployeeDetails.xlsx - "UI Fixes and some backend fix fo import"
#
# 1) Record a new employee row (row 95) that was imported into the
#    "Jul 2022" sheet - Employee #189 / "Test Import" / "abc" job title &
#    department - right above the generated "Generated on ..." footer row.
# 2) Leave the viewport scrolled down to / focused on the newly added row,
#    matching where the user was working when they saved.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jul 2022")

# --- Add the new employee record in row 95 ---------------------------------
# Copy the formatting from the row above first so the date columns (C/F)
# keep the same date number-format style as the rest of the table.
$ws.Range("A94:F94").Copy()
$ws.Range("A95:F95").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(95, 1).Value = 189
$ws.Cells.Item(95, 2).Value = "Test Import"
$ws.Cells.Item(95, 3).Value = 44753
$ws.Cells.Item(95, 4).Value = "abc"
$ws.Cells.Item(95, 5).Value = "abc"
$ws.Cells.Item(95, 6).Value = 35800

$excel.CutCopyMode = $false

# --- Restore the view / selection so the sheet reopens scrolled to the
#     newly imported row, matching the author's saved window state. -------
$win = $excel.ActiveWindow
$win.ScrollRow = 85
$win.ScrollColumn = 1
$win.Left = -108
$win.Top = -108
$win.Width = 23256
$win.Height = 12456

[void]$ws.Range("C95").Select()

Write-Host "Inserted row 95 (Employee 189 / Test Import) and updated the view."
